$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-02-27 Friday" "2026-02-28 Saturday"

Replace-Text "571×9=5139" "865×2=1730"
Replace-Text "467×9=4203" "454×5=2270"
Replace-Text "591×7=4137" "570×9=5130"
Replace-Text "596×8=4768" "227×9=2043"
Replace-Text "755×6=4530" "698×3=2094"

Replace-Text "726×6=4356" "656×3=1968"
Replace-Text "129×3=387" "301×4=1204"
Replace-Text "956×2=1912" "443×9=3987"
Replace-Text "175×5=875" "955×9=8595"
Replace-Text "534×2=1068" "450×6=2700"

Replace-Text "281×2=562" "459×7=3213"
Replace-Text "606×8=4848" "166×7=1162"
Replace-Text "570×2=1140" "139×6=834"
Replace-Text "327×9=2943" "706×3=2118"
Replace-Text "165×4=660" "974×8=7792"

Replace-Text "515×5=2575" "559×2=1118"
Replace-Text "839×9=7551" "545×7=3815"
Replace-Text "904×6=5424" "251×5=1255"
Replace-Text "520×2=1040" "251×3=753"
Replace-Text "306×2=612" "239×8=1912"

Replace-Text "234×3=702" "391×6=2346"
Replace-Text "919×5=4595" "491×6=2946"
Replace-Text "543×3=1629" "814×6=4884"
Replace-Text "895×7=6265" "778×9=7002"
Replace-Text "255×2=510" "910×7=6370"
